$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text string, forcing text storage even when
# the string looks numeric (matches the source data which stores Price/Volume
# columns as plain text, e.g. "253.05", "  -1.59%  ").
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    if ($value -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $value
}

# Full row replacements (coin reordering: rows 9/10 swap USDC<->Cardano, rows 31/32 swap Dai<->Hedera)
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D9" "1.05"
Set-TextValue "E9" "  -1.33%  "

$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D10" "0.999"
Set-TextValue "E10" "  +0.02%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D31" "0.142"
Set-TextValue "E31" "  -5.92%  "

$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  +0.63%  "

# Price (D) and Volume(1h) (E) updates
Set-TextValue "D2" "97.910.36"
Set-TextValue "E2" "  -0.35%  "
Set-TextValue "D3" "3.378.89"
Set-TextValue "E3" "  -1.11%  "
Set-TextValue "D5" "253.05"
Set-TextValue "E5" "  -1.59%  "
Set-TextValue "D6" "661.68"
Set-TextValue "E6" "  +0.19%  "
Set-TextValue "D7" "1.46"
Set-TextValue "E7" "  -0.56%  "
Set-TextValue "D8" "0.423"
Set-TextValue "E8" "  -3.49%  "
Set-TextValue "D11" "3.375.55"
Set-TextValue "E11" "  -1.11%  "
Set-TextValue "D13" "41.53"
Set-TextValue "E13" "  -2.56%  "
Set-TextValue "D14" "97.471.98"
Set-TextValue "E14" "  -0.44%  "
Set-TextValue "D15" "6.11"
Set-TextValue "E15" "  -4.76%  "
Set-TextValue "D16" "0.0000255"
Set-TextValue "E16" "  -4.77%  "
Set-TextValue "D17" "4.011.47"
Set-TextValue "E17" "  -1.11%  "
Set-TextValue "D18" "8.90"
Set-TextValue "E18" "  -3.01%  "
Set-TextValue "D19" "3.368.17"
Set-TextValue "E19" "  -1.61%  "
Set-TextValue "D20" "18.07"
Set-TextValue "E20" "  +1.57%  "
Set-TextValue "D21" "0.533"
Set-TextValue "E21" "  -10.46%  "
Set-TextValue "D22" "10.92"
Set-TextValue "E22" "  -0.45%  "
Set-TextValue "D23" "510.74"
Set-TextValue "E23" "  -1.12%  "
Set-TextValue "D25" "7.03"
Set-TextValue "E25" "  +8.54%  "
Set-TextValue "D26" "0.0000200"
Set-TextValue "E26" "  -3.40%  "
Set-TextValue "D27" "96.59"
Set-TextValue "E27" "  -4.53%  "
Set-TextValue "D28" "12.34"
Set-TextValue "E28" "  -5.35%  "
Set-TextValue "D29" "3.554.39"
Set-TextValue "E29" "  -1.19%  "
Set-TextValue "D30" "11.41"
Set-TextValue "E30" "  -3.51%  "
Set-TextValue "D33" "0.187"
Set-TextValue "E33" "  -5.56%  "
Set-TextValue "D34" "2.58"
Set-TextValue "E34" "  +8.05%  "
Set-TextValue "D35" "0.998"
Set-TextValue "E35" "  -0.03%  "
Set-TextValue "D36" "0.561"
Set-TextValue "E36" "  -3.91%  "
Set-TextValue "D37" "28.79"
Set-TextValue "E37" "  -4.52%  "
Set-TextValue "D38" "7.98"
Set-TextValue "E38" "  +0.53%  "
Set-TextValue "D39" "1.50"
Set-TextValue "E39" "  +0.01%  "
Set-TextValue "D40" "531.10"
Set-TextValue "E40" "  -0.99%  "
Set-TextValue "D41" "0.152"
Set-TextValue "E41" "  -0.95%  "
Set-TextValue "D43" "24.40"
Set-TextValue "E43" "  -1.25%  "
Set-TextValue "D44" "0.854"
Set-TextValue "E44" "  -3.49%  "
Set-TextValue "D48" "3.67"
Set-TextValue "E48" "  -1.45%  "
Set-TextValue "D49" "5.63"
Set-TextValue "E49" "  -5.11%  "
Set-TextValue "D50" "56.16"
Set-TextValue "E50" "  +2.17%  "
Set-TextValue "D51" "8.62"
Set-TextValue "E51" "  -6.50%  "

# Volume(1h) (E) only updates
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "E12" "  -2.87%  "
Set-TextValue "E24" "  -1.77%  "
Set-TextValue "E42" "  +0.00%  "
Set-TextValue "E45" "  -0.63%  "
Set-TextValue "E46" "  +1.97%  "
Set-TextValue "E47" "  +8.35%  "
